# Updates cryptos list figures (prices / 1h volume %) to match latest
# scrape, per commit "Updated cryptos list ... with GitHub Actions".
# All target cells are plain text in the sheet (inline strings), so
# numeric-looking price values are written with a leading quote
# (Excel's text quote-prefix) to keep them stored as text instead of
# being auto-coerced to numbers (e.g. "100.90" -> 100.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.006.26'
$ws.Range('E2').Value = '  +7.91%  '
$ws.Range('D3').Value = '1.821.27'
$ws.Range('E3').Value = '  +5.25%  '
$ws.Range('D4').Value = '''0.9998'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''246.28'
$ws.Range('E5').Value = '  +2.63%  '
$ws.Range('D6').Value = '''0.9997'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = '''0.4931'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('D8').Value = '''44.34'
$ws.Range('E8').Value = '  +7.16%  '
$ws.Range('E9').Value = '  +6.44%  '
$ws.Range('D10').Value = '''0.06381'
$ws.Range('E10').Value = '  +3.31%  '
$ws.Range('D11').Value = '1.819.95'
$ws.Range('E11').Value = '  +5.21%  '
$ws.Range('E12').Value = '  +3.69%  '
$ws.Range('D13').Value = '''0.07053'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').Value = '''0.6431'
$ws.Range('D15').Value = '''84.01'
$ws.Range('E15').Value = '  +9.14%  '
$ws.Range('D16').Value = '''4.692'
$ws.Range('E16').Value = '  +5.14%  '
$ws.Range('D17').Value = '29.012.85'
$ws.Range('E17').Value = '  +8.86%  '
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '''0.000007288'
$ws.Range('E19').Value = '  +2.28%  '
$ws.Range('B20').Value = 'BinanceUSD'
$ws.Range('C20').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D20').Value = '''0.9996'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '''12.19'
$ws.Range('E21').Value = '  +7.22%  '
$ws.Range('D22').Value = '2.055.20'
$ws.Range('E22').Value = '  +5.33%  '
$ws.Range('D23').Value = '''4.547'
$ws.Range('E23').Value = '  +3.43%  '
$ws.Range('D24').Value = '''8.832'
$ws.Range('E24').Value = '  +4.89%  '
$ws.Range('D25').Value = '''5.364'
$ws.Range('E25').Value = '  +5.97%  '
$ws.Range('D26').Value = '''143.84'
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('D27').Value = '''129.48'
$ws.Range('E27').Value = '  +21.53%  '
$ws.Range('D28').Value = '''16.36'
$ws.Range('E28').Value = '  +7.63%  '
$ws.Range('D29').Value = '''1.881'
$ws.Range('E29').Value = '  +4.68%  '
$ws.Range('D30').Value = '''1.399'
$ws.Range('E30').Value = '  +1.72%  '
$ws.Range('D31').Value = '''4.123'
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('D32').Value = '''0.08354'
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('D33').Value = '''3.770'
$ws.Range('E33').Value = '  +2.87%  '
$ws.Range('D34').Value = '''0.04952'
$ws.Range('E34').Value = '  +8.20%  '
$ws.Range('D35').Value = '''1.097'
$ws.Range('E35').Value = '  +9.65%  '
$ws.Range('D36').Value = '''2.706'
$ws.Range('E36').Value = '  +4.37%  '
$ws.Range('D37').Value = '''0.6687'
$ws.Range('E37').Value = '  +8.43%  '
$ws.Range('D38').Value = '''2.295'
$ws.Range('E38').Value = '  +15.63%  '
$ws.Range('D39').Value = '''2.684'
$ws.Range('E39').Value = '  +9.16%  '
$ws.Range('D40').Value = '''0.9467'
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('D41').Value = '''6.145'
$ws.Range('E41').Value = '  +8.16%  '
$ws.Range('D42').Value = '''0.01582'
$ws.Range('E42').Value = '  +5.70%  '
$ws.Range('D43').Value = '''0.9999'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = '''100.90'
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').Value = '''0.4059'
$ws.Range('E45').Value = '  +5.77%  '
$ws.Range('D46').Value = '''7.169'
$ws.Range('E46').Value = '  +5.70%  '
$ws.Range('D47').Value = '''0.1219'
$ws.Range('E47').Value = '  +5.63%  '
$ws.Range('D48').Value = '''0.05545'
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('D49').Value = '''31.63'
$ws.Range('E49').Value = '  +5.29%  '
$ws.Range('D50').Value = '''8.091'
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('E51').Value = '  +4.56%  '
